$wb = $excel.ActiveWorkbook

# Reference sheets already present
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws4 = $wb.Worksheets.Item($wb.Worksheets.Count)

# --- Add sheet 5: getAllCacheNames ---
$ws5 = $wb.Worksheets.Add($null, $ws4)
$ws5.Name = "getAllCacheNames"

# Bring over the existing header look (font/fill/border) from an existing header row,
# then drop the vertical-center alignment that the new header style doesn't use.
$null = $ws1.Range("A1:F1").Copy($ws5.Range("A1:F1"))
$ws5.Range("A1:F1").VerticalAlignment = -4107

$ws5.Range("A1").Value = "test-id"
$ws5.Range("B1").Value = "description"
$ws5.Range("C1").Value = "dataList"
$ws5.Range("D1").Value = "rspStatus"
$ws5.Range("E1").Value = "rspCode"
$ws5.Range("F1").Value = "rspMessage"

$ws5.Range("A2").Value = "jinzu-connector-configure-cache-test1"
$ws5.Range("B2").Value = "get allCacheNames"
$ws5.Range("C2").Value = "auth,connector,dataSource,engine,entity,mapper,plugin,rule,transaction,licenseConfigEncrypted"

$ws5.Columns.Item(1).ColumnWidth = 39
$ws5.Columns.Item(2).ColumnWidth = 26
$ws5.Columns.Item(3).ColumnWidth = 26

$null = $ws5.Range("B10").Select()

# --- Add sheet 6: getCacheKeyAndValue ---
$ws6 = $wb.Worksheets.Add($null, $ws5)
$ws6.Name = "getCacheKeyAndValue"

$null = $ws2.Range("A1:G1").Copy($ws6.Range("A1:G1"))
$ws6.Range("A1:G1").VerticalAlignment = -4107

$ws6.Range("A1").Value = "test-id"
$ws6.Range("B1").Value = "description"
$ws6.Range("C1").Value = "cacheName"
$ws6.Range("D1").Value = "entityName"
$ws6.Range("E1").Value = "rspStatus"
$ws6.Range("F1").Value = "rspCode"
$ws6.Range("G1").Value = "rspMessage"

$ws6.Range("A2").Value = "jinzu-connector-configure-cache-test2"
$ws6.Range("B2").Value = " mapper,check get cache key and value"
$ws6.Range("C2").Value = "mapper"
$ws6.Range("D2").Value = "Site"

$ws6.Columns.Item(1).ColumnWidth = 27
$ws6.Columns.Item(2).ColumnWidth = 42
$ws6.Columns.Item(3).ColumnWidth = 22
$ws6.Columns.Item(4).ColumnWidth = 30

$null = $ws6.Range("B11").Select()
$null = $ws6.Activate()
